$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Botoks"
$ws.Range("B4").Value = "Zabieg et"
$ws.Range("C4").Value = 1200

$ws.Range("A5").Value = "Oczyszczanie"
$ws.Range("B5").Value = "Opis"
$ws.Range("C5").Value = 500
